$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values that changed
$ws.Range("B286").Value = 197.001
$ws.Range("B293").Value = 205.5844
$ws.Range("B301").Value = 189.7653

# New rows with dates and values
$dates = @("28-10-2021", "29-10-2021", "30-10-2021", "31-10-2021", "01-11-2021", "02-11-2021")
$values = @(183.1283, 179.3698, 174.2133, 166.1266, 167.3835, 173.1238)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 302 + $i
    $cellA = $ws.Cells.Item($row, 1)
    # Force text so date-looking strings like "01-11-2021" are not
    # auto-converted into date serial numbers by Excel.
    $cellA.NumberFormat = "@"
    $cellA.Value = $dates[$i]
    $cellA.ClearFormats()
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
